$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Update rule name in row 11 (was "R40") to the text value "1".
# Leading apostrophe forces Excel to store it as text (shared string)
# instead of auto-converting it to a numeric literal.
$ws.Range("B11").Value = "'1"
